$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every Price cell (column D) as literal text,
# including values that look numeric (e.g. "0.998", "1.00", "37.87").
# A plain numeric-looking string assigned via .Value gets auto-converted
# by Excel into a real number, which would lose formatting (trailing
# zeros, etc). For those cells we briefly force Text number format so
# the literal text is kept, then restore the "Normal" style so the
# cell format matches the rest of the (unstyled) data rows.

$ws.Range("D2").Value = "66.497.92"
$ws.Range("E2").Value = "  -5.12%  "
$ws.Range("D3").Value = "3.218.63"
$ws.Range("E3").Value = "  -8.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.94%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.206.71"
$ws.Range("E8").Value = "  -8.39%  "
$ws.Range("E9").Value = "  -13.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -19.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -17.89%  "
$ws.Range("E14").Value = "  -15.30%  "
$ws.Range("D15").Value = "3.729.48"
$ws.Range("E15").Value = "  -8.38%  "
$ws.Range("D16").Value = "66.286.01"
$ws.Range("E16").Value = "  -5.55%  "
$ws.Range("D17").Value = "3.203.23"
$ws.Range("E17").Value = "  -8.80%  "
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "518.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -15.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -18.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -16.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.744"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -14.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -16.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -14.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -15.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -13.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "28.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -15.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -18.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -17.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -15.40%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "530.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -16.69%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -17.23%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -21.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0416"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -16.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -15.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.76%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -23.51%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.827.51"
$ws.Range("E43").Value = "  -15.60%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -17.09%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.249"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -18.68%  "
$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").Value = "0.0₃0548"
$ws.Range("E47").Value = "  -25.37%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -20.38%  "
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -18.57%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.57%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -14.16%  "
